# Wyoming 2023 MCAS workbook cleanup:
#  1. Rename header row to snake_case field names.
#  2. Convert the all-caps state/municipality names (columns A & B,
#     rows 2-347) to Title Case.
#  3. Remove the trailing metadata/footer rows (349-353) that are no
#     longer part of the clean dataset. This also shrinks the sheet
#     dimension from A1:D353 down to A1:D347 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns (row 1) ---
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# --- 2. Title-case the state (col A) and municipality (col B) names ---
$textInfo = (Get-Culture).TextInfo
for ($r = 2; $r -le 347; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null) {
            $cell.Value = $textInfo.ToTitleCase($v.ToLower())
        }
    }
}

# --- 3. Delete the footer/metadata rows 349-353 ---
$ws.Range("A349:A353").EntireRow.Delete()
